# Correct preparer name "H.Brown" -> "H.BROWN" for the s2cDNAPreparer (B)
# and libraryPreparer (E) columns, and apply an explicit black font color
# to those cells (matches the standardized formatting used elsewhere).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace all occurrences of "H.Brown" with "H.BROWN" across the sheet.
[void]$ws.Cells.Replace("H.Brown", "H.BROWN")

# Apply explicit black font color to the preparer columns (B and E, rows 2-27).
$ws.Range("B2:B27").Font.Color = 0
$ws.Range("E2:E27").Font.Color = 0

# Reflect the selection left after the edit (column E, the last column
# touched interactively by the author).
[void]$ws.Range("E2:E27").Select()
